$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9528933763504028
$ws.Range("B1").Value = 1.770101189613342
$ws.Range("C1").Value = 4.974975109100342
$ws.Range("D1").Value = 2.182043552398682
$ws.Range("E1").Value = 1.369734287261963
